# Completed some of the simple tasks in CodeChanges.xlsx except adding the
# try catch block.
#
# Mark the "Done" / "Fuck you" status for each completed task in column B
# (rows with "Use try catch block" in column E - rows 5 & 8 - are left
# untouched, since that task was not completed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = "Done"
$ws.Range("B4").Value  = "Done"
$ws.Range("B6").Value  = "Done"
$ws.Range("B7").Value  = "Fuck you"
$ws.Range("B9").Value  = "Done"
$ws.Range("B10").Value = "Done"

# Set up the page for printing (A4, portrait).
$pageSetup = $ws.PageSetup
$pageSetup.PaperSize = 9
$pageSetup.Orientation = 1

# Leave the selection on B8, matching where the author ended up working.
$ws.Range("B8").Select() | Out-Null
